$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 (Caulfield train exposure site)
$ws.Range("A2").Value = "Caulfield"
$ws.Range("B2").Value = "Metro Train - Frankston line"
$ws.Range("C2").Value = "30/12/20 4:30pm-17:00pm"
$ws.Range("D2").Value = "Case caught train from Caulfield to Cheltenham"
$ws.Range("E2").Value = "new"

# Update row 3 (Lakes Entrance bus exposure site)
$ws.Range("A3").Value = "Lakes Entrance"
$ws.Range("B3").Value = "V/Line bus - Lakes Entrance to Bairnsdale"
$ws.Range("C3").Value = "30/12/2020 11:55am-12:30pm"
$ws.Range("D3").Value = "Case caught the 11:55am bus from Lakes Entrance"
$ws.Range("E3").Value = "new"

# Remove rows 4-6 which held the old extra exposure sites
$ws.Range("A4:E6").Delete()

# Re-fit column widths to new (shorter) content, matching Excel's recalculated
# best-fit widths for the new text
$ws.Columns.Item(2).ColumnWidth = 32.5
$ws.Columns.Item(3).ColumnWidth = 24.833333333333336
$ws.Columns.Item(4).ColumnWidth = 39.83333333333333

# Update selection to whole columns A:E as in final file
$ws.Range("A1:E1048576").Select()
